$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update checklist status values (accessToken implemented, refreshToken almost done)
$ws.Range("C8").Value = "ok"
$ws.Range("C9").Value = "ok"
$ws.Range("C10").Value = "Em andamento"

# Update the active selection on the sheet
[void]$ws.Range("D8").Select()
